# B6-PowerPoint.pptx edit:
#  1. Re-style the three "Component three" tables (slides 14, 15, 16) from the
#     bespoke "Table_0" style to the built-in "Medium Style 2 - Accent 2" table
#     style ({B0F05431-4BF8-4B80-BEE0-A6B2660F8B81}), as happens when a user
#     picks a new look from the Table Design gallery.
#  2. Re-apply the deck's colour scheme (this presentation's single Design /
#     slide master) so its RGB values match the plain "Office Theme" palette
#     instead of the "Integral / Red Violet" palette it currently carries.

$p = $ppt.ActivePresentation

# --- 1. Table style swap -----------------------------------------------
$newStyleId = "{B0F05431-4BF8-4B80-BEE0-A6B2660F8B81}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme swap ----------------------------------------
# The presentation's theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink - in that order) are driven back to the stock "Office Theme"
# values. Colours are COM `RGB` longs, i.e. 0x00BBGGRR.

function ToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 0; $i -lt $officeThemeHex.Count; $i++) {
    $tcs.Item($i + 1).RGB = ToComRGB $officeThemeHex[$i]
}
